# average with safety stocks
$wb = $excel.ActiveWorkbook

# --- Productdata sheet: SetupCosts (D), AverageDemand (F), StandardDevDemands (I)
# scaled down by a factor of 0.0004 for rows 2-11 ---
$ws = $wb.Worksheets.Item("Productdata")

$setupCosts = @{2=0.0016; 3=0.0028; 4=0.0024; 5=0.0012; 6=0.0012; 7=0.0012; 8=0.0008; 9=0.0004; 10=0.0004; 11=0.0004}
$avgDemand  = @{2=0.0032; 3=0.0056; 4=0.0048; 5=0.0024; 6=0.0024; 7=0.0024; 8=0.0016; 9=0.0008; 10=0.0008; 11=0.0008}
$stdDevDem  = @{2=0.032;  3=0.056;  4=0.048;  5=0.024;  6=0.024;  7=0.024;  8=0.016;  9=0.008;  10=0.008;  11=0.008}

foreach ($row in 2..11) {
    $ws.Cells.Item($row, 4).Value = $setupCosts[$row]
    $ws.Cells.Item($row, 6).Value = $avgDemand[$row]
    $ws.Cells.Item($row, 9).Value = $stdDevDem[$row]
}

# --- ForcastedStandardDeviation sheet: zero out rows 9-11, columns B-E ---
$ws2 = $wb.Worksheets.Item("ForcastedStandardDeviation")

foreach ($row in 9..11) {
    foreach ($col in 2..5) {
        $ws2.Cells.Item($row, $col).Value = 0
    }
}
